$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write new vocabulary rows (N5/N4/N3 additions) ---
$ws.Range('B85').Value = 'かばん'
$ws.Range('C85').Value = ' 皮包'
$ws.Rows.Item(85).RowHeight = 18.75

$ws.Range('A86').Value = '花瓶'
$ws.Range('B86').Value = 'かびん'
$ws.Range('C86').Value = '花瓶'
$ws.Rows.Item(86).RowHeight = 18.75

$ws.Range('A87').Value = '紙'
$ws.Range('B87').Value = 'かみ'
$ws.Range('C87').Value = '紙'
$ws.Rows.Item(87).RowHeight = 18.75

$ws.Range('A88').Value = '体'
$ws.Range('B88').Value = 'からだ'
$ws.Range('C88').Value = '身體'
$ws.Rows.Item(88).RowHeight = 18.75

$ws.Range('A89').Value = '川'
$ws.Range('B89').Value = 'かわ'
$ws.Range('C89').Value = '河川'
$ws.Rows.Item(89).RowHeight = 18.75

$ws.Range('A90').Value = '漢字'
$ws.Range('B90').Value = 'かんじ'
$ws.Range('C90').Value = '漢字'
$ws.Rows.Item(90).RowHeight = 18.75

$ws.Range('A91').Value = '木'
$ws.Range('B91').Value = 'き'
$ws.Range('C91').Value = '樹'
$ws.Rows.Item(91).RowHeight = 18.75

$ws.Range('A92').Value = '北'
$ws.Range('B92').Value = 'きた'
$ws.Range('C92').Value = '北方北邊'
$ws.Rows.Item(92).RowHeight = 18.75

$ws.Range('A93').Value = '喫茶店'
$ws.Range('B93').Value = 'きっさてん'
$ws.Range('C93').Value = '咖啡廳'
$ws.Rows.Item(93).RowHeight = 18.75

$ws.Range('A94').Value = '切手'
$ws.Range('B94').Value = 'きって'
$ws.Range('C94').Value = '郵票'
$ws.Rows.Item(94).RowHeight = 18.75

$ws.Range('A95').Value = '切符'
$ws.Range('B95').Value = 'きっぷ'
$ws.Range('C95').Value = '車票'
$ws.Rows.Item(95).RowHeight = 18.75

$ws.Range('A96').Value = '昨日'
$ws.Range('B96').Value = 'きのう'
$ws.Range('C96').Value = '昨天'
$ws.Rows.Item(96).RowHeight = 18.75

$ws.Range('A97').Value = '牛肉'
$ws.Range('B97').Value = 'ぎゅうにく'
$ws.Range('C97').Value = '牛肉'
$ws.Rows.Item(97).RowHeight = 18.75

$ws.Range('A98').Value = '牛乳'
$ws.Range('B98').Value = 'ぎゅうにゆう'
$ws.Range('C98').Value = '牛奶'
$ws.Rows.Item(98).RowHeight = 18.75

$ws.Range('A99').Value = '今日'
$ws.Range('B99').Value = 'きょう'
$ws.Range('C99').Value = '今天'
$ws.Rows.Item(99).RowHeight = 18.75

$ws.Range('A100').Value = '教室'
$ws.Range('B100').Value = 'きょうしつ'
$ws.Range('C100').Value = '教室'
$ws.Rows.Item(100).RowHeight = 18.75

$ws.Range('A101').Value = '兄弟'
$ws.Range('B101').Value = 'きょうだい'
$ws.Range('C101').Value = '兄弟'
$ws.Rows.Item(101).RowHeight = 18.75

$ws.Range('A102').Value = '去年'
$ws.Range('B102').Value = 'きょねん'
$ws.Range('C102').Value = '去年'
$ws.Rows.Item(102).RowHeight = 18.75

$ws.Range('A103').Value = '銀行'
$ws.Range('B103').Value = 'ぎんこう'
$ws.Range('C103').Value = '銀行'
$ws.Rows.Item(103).RowHeight = 18.75

$ws.Range('A104').Value = '薬'
$ws.Range('B104').Value = 'くすり'
$ws.Range('C104').Value = '藥'
$ws.Rows.Item(104).RowHeight = 18.75

$ws.Range('A105').Value = '果物'
$ws.Range('B105').Value = 'くだもの'
$ws.Range('C105').Value = '水果'
$ws.Rows.Item(105).RowHeight = 18.75

$ws.Range('A106').Value = '口'
$ws.Range('B106').Value = 'くち'
$ws.Range('C106').Value = '嘴巴'
$ws.Rows.Item(106).RowHeight = 18.75

$ws.Range('A107').Value = '靴'
$ws.Range('B107').Value = 'くつ'
$ws.Range('C107').Value = '鞋子'
$ws.Rows.Item(107).RowHeight = 18.75

$ws.Range('A108').Value = '靴下'
$ws.Range('B108').Value = 'くつした'
$ws.Range('C108').Value = '襪子'
$ws.Rows.Item(108).RowHeight = 18.75

$ws.Range('A109').Value = '国'
$ws.Range('B109').Value = 'くに'
$ws.Range('C109').Value = '國家'
$ws.Rows.Item(109).RowHeight = 18.75

$ws.Range('A110').Value = '車'
$ws.Range('B110').Value = 'くるま'
$ws.Range('C110').Value = '車子'
$ws.Rows.Item(110).RowHeight = 18.75

$ws.Range('A111').Value = '警官'
$ws.Range('B111').Value = 'けいかん'
$ws.Range('C111').Value = '警察'
$ws.Rows.Item(111).RowHeight = 18.75

$ws.Range('A112').Value = '今朝'
$ws.Range('B112').Value = 'けさ'
$ws.Range('C112').Value = '今天早上'
$ws.Rows.Item(112).RowHeight = 18.75

$ws.Range('A113').Value = '結婚'
$ws.Range('B113').Value = 'けっこん'
$ws.Range('C113').Value = '結婚'
$ws.Rows.Item(113).RowHeight = 18.75

$ws.Range('A114').Value = '公園'
$ws.Range('B114').Value = 'こうえん'
$ws.Range('C114').Value = '公園'
$ws.Rows.Item(114).RowHeight = 18.75

$ws.Range('A115').Value = '紅茶'
$ws.Range('B115').Value = 'こうちゃ'
$ws.Range('C115').Value = '紅茶'
$ws.Rows.Item(115).RowHeight = 18.75

$ws.Range('A116').Value = '交番'
$ws.Range('B116').Value = 'こうばん'
$ws.Range('C116').Value = '警察局'
$ws.Rows.Item(116).RowHeight = 18.75

$ws.Range('A117').Value = '声'
$ws.Range('B117').Value = 'こえ'
$ws.Range('C117').Value = '聲音'
$ws.Rows.Item(117).RowHeight = 18.75

$ws.Range('A118').Value = '午後'
$ws.Range('B118').Value = 'ごご'
$ws.Range('C118').Value = '下午'
$ws.Rows.Item(118).RowHeight = 18.75

$ws.Range('A119').Value = '午前'
$ws.Range('B119').Value = 'ごぜん'
$ws.Range('C119').Value = '上午'
$ws.Rows.Item(119).RowHeight = 18.75

$ws.Range('B120').Value = 'こちら'
$ws.Range('C120').Value = '(尊敬)這邊'
$ws.Rows.Item(120).RowHeight = 18.75

$ws.Range('A121').Value = '今年'
$ws.Range('B121').Value = 'ことし'
$ws.Range('C121').Value = '今年'
$ws.Rows.Item(121).RowHeight = 18.75

$ws.Range('A122').Value = '子ども'
$ws.Range('B122').Value = 'こども'
$ws.Range('C122').Value = '小孩子'
$ws.Rows.Item(122).RowHeight = 18.75

$ws.Range('A123').Value = 'ご飯'
$ws.Range('B123').Value = 'ごはん'
$ws.Range('C123').Value = '吃飯，飯'
$ws.Rows.Item(123).RowHeight = 18.75

$ws.Range('A124').Value = '今月'
$ws.Range('B124').Value = 'こんげつ'
$ws.Range('C124').Value = '這個月'
$ws.Rows.Item(124).RowHeight = 18.75

$ws.Range('A125').Value = '今週'
$ws.Range('B125').Value = 'こんしゅう'
$ws.Range('C125').Value = '這星期'
$ws.Rows.Item(125).RowHeight = 18.75

$ws.Range('A126').Value = '今晩'
$ws.Range('B126').Value = 'こんばん'
$ws.Range('C126').Value = '今晚'
$ws.Rows.Item(126).RowHeight = 18.75

# --- Apply Yu Gothic font to Japanese word/reading columns (matches original author formatting) ---
$ws.Range("B85:B126").Font.Name = "Yu Gothic"
$ws.Range("A87:A119").Font.Name = "Yu Gothic"
$ws.Range("A121:A126").Font.Name = "Yu Gothic"

# --- Page setup (print: A4/Letter #9, portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore view/selection state ---
$ws.Application.Goto($ws.Range("A113"), $true)
$ws.Range("F120").Select()